$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("Input")
$wsOutput = $wb.Worksheets.Item("Output")

# --- Input sheet content changes ---
$wsInput.Range("B6").Value = "Jhon Deer"
$wsInput.Range("A7").Value = "GroupAddClient"
$wsInput.Range("B7").Value = "click"

# --- Output sheet content changes ---
# A1 now carries the same "label" formatting used by column A on the Input
# sheet (Calibri 11 on the grey header fill), so copy that formatting over
# before changing the text.
$wsInput.Range("A1").Copy()
$wsOutput.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

$wsOutput.Range("A1").Value = "verify1"
$wsOutput.Range("B1").Value = "Cannot close. Group has active clients"

# Column B has to widen considerably to fit the new, much longer message
$wsOutput.Columns.Item(2).ColumnWidth = 32.6

# Row 1 grows to match the taller font now used in A1
$wsOutput.Rows.Item(1).RowHeight = 15

# --- Selections / active sheet ---
# Input is no longer the active/selected tab; its selection moves to C17
$wsInput.Range("C17").Select()

# Output becomes the active tab with its selection on B5
$wsOutput.Activate()
$wsOutput.Range("B5").Select()
